$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concurrent Iterations")
$ws.Activate()

$ws.Range("A22").Value = "Determining salary at any given year"
$ws.Range("A23").Value = "Determining wealth of an individual after he saves"
$ws.Range("A24").Value = "Determining years to retirement based on individual's wealth"

$ws.Range("A22:A24").Select()
